$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knowledge Based Agent")

$data = @(
    @(4, 1, 0, 0, 0, 4, 996),
    @(7, 1, 0, 0, 0, 7, 993),
    @(2, 1, 0, 0, 0, 2, 998),
    @(15, 0, 0, 0, 0, 14, -22),
    @(16, 0, 0, 0, 0, 15, -22),
    @(3, 0, 0, 0, 0, 2, -3),
    @(20, 0, 0, 0, 0, 19, -28),
    @(5, 1, 0, 0, 0, 5, 994),
    @(6, 0, 0, 0, 0, 5, -8),
    @(5, 1, 0, 0, 0, 5, 995),
    @(12, 1, 0, 0, 0, 12, 985),
    @(19, 0, 0, 0, 0, 18, -24),
    @(17, 1, 0, 0, 0, 17, 978),
    @(14, 1, 0, 0, 0, 14, 979),
    @(12, 1, 0, 0, 0, 12, 982),
    @(3, 0, 0, 0, 0, 2, -3),
    @(4, 1, 0, 0, 0, 4, 996),
    @(7, 1, 0, 0, 0, 7, 988),
    @(3, 0, 0, 0, 0, 2, -3),
    @(13, 0, 0, 0, 0, 12, -19)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($row, 2 + $col).Value = $values[$col]
    }
}

$ws.Range("L16").Select()
